$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values - keep as text to match source formatting
$ws.Range("D2").Value = "'246.12"
$ws.Range("D3").Value = "'22.77"
$ws.Range("D4").Value = "'5.277"
$ws.Range("D7").Value = "'0.8103"
$ws.Range("D8").Value = "'0.8823"
$ws.Range("D9").Value = "'0.1425"
$ws.Range("D10").Value = "'0.07368"
$ws.Range("D11").Value = "'0.03002"
$ws.Range("D12").Value = "'0.03128"
$ws.Range("D13").Value = "'0.09396"
$ws.Range("D14").Value = "'3.928"
$ws.Range("D15").Value = "'0.001578"
$ws.Range("D17").Value = "'0.0005844"
$ws.Range("D18").Value = "'0.006141"
$ws.Range("D19").Value = "'0.005096"
$ws.Range("D20").Value = "'0.0009986"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("D22").Value = "'3.745"
$ws.Range("D23").Value = "'6.306"
$ws.Range("D24").Value = "'2.186"
$ws.Range("D25").Value = "'0.3279"
$ws.Range("D27").Value = "'0.0001101"
$ws.Range("D40").Value = "'0.03906"
$ws.Range("D41").Value = "'0.006735"
$ws.Range("D42").Value = "'0.1065"
$ws.Range("D43").Value = "'0.003202"
$ws.Range("D44").Value = "'0.007503"
$ws.Range("D45").Value = "'0.00005637"
$ws.Range("D47").Value = "'0.6004"
$ws.Range("D48").Value = "'0.1744"
$ws.Range("D50").Value = "'0.01011"

# Update Volume(1h) text (column E) for rows 41 and 47
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# Update Hora (column G) from 6 to 7 for all data rows (2-51)
$ws.Range("G2:G51").Value = "'7"

